# Images and PCB Update
# - Rename sheet tab to reflect new PCB export date
# - Update existing U1 (formerly mis-set as "U2") row with corrected pin count / SMD flag / Pad X
# - Append newly placed components (C1, C2, R1, U2) pick-and-place rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to the new PCB export date
$ws.Name = "PickAndPlace_PCB1_2025-11-27"

# --- Row 3: fix up the ESP32-S3-Zero row (designator, pad X, pin count, SMD flag) ---
$ws.Cells.Item(3, 1).Value  = "U1"                                       # Designator
$ws.Cells.Item(3, 2).Value  = "ESP32-S3-Zero"                            # Device
$ws.Cells.Item(3, 3).Value  = "COMM-SMD_18P-P2.54-L23.5-W18.0-TL"        # Footprint
$ws.Cells.Item(3, 4).Value  = "15.367mm"                                 # Mid X
$ws.Cells.Item(3, 5).Value  = "-14.859mm"                                # Mid Y
$ws.Cells.Item(3, 6).Value  = "15.367mm"                                 # Ref X
$ws.Cells.Item(3, 7).Value  = "-14.859mm"                                # Ref Y
$ws.Cells.Item(3, 8).Value  = "6.367mm"                                  # Pad X
$ws.Cells.Item(3, 9).Value  = "-4.699mm"                                 # Pad Y
$ws.Cells.Item(3, 10).Value = 18                                         # Pins
$ws.Cells.Item(3, 11).Value = "T"                                        # Layer
$ws.Cells.Item(3, 12).Value = 0                                          # Rotation
$ws.Cells.Item(3, 13).Value = "Yes"                                      # SMD
$ws.Cells.Item(3, 14).Value = "ESP32-S3-Zero"                            # Comment

# --- Row 4: C1 capacitor ---
$ws.Cells.Item(4, 1).Value  = "C1"
$ws.Cells.Item(4, 2).Value  = "CL05A105KO5NNNC"
$ws.Cells.Item(4, 3).Value  = "C0402"
$ws.Cells.Item(4, 4).Value  = "25.4mm"
$ws.Cells.Item(4, 5).Value  = "-29.845mm"
$ws.Cells.Item(4, 6).Value  = "25.4mm"
$ws.Cells.Item(4, 7).Value  = "-29.845mm"
$ws.Cells.Item(4, 8).Value  = "25.4mm"
$ws.Cells.Item(4, 9).Value  = "-30.39mm"
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 11).Value = "T"
$ws.Cells.Item(4, 12).Value = 90
$ws.Cells.Item(4, 13).Value = "Yes"
$ws.Cells.Item(4, 14).Value = "1uF"

# --- Row 5: C2 capacitor ---
$ws.Cells.Item(5, 1).Value  = "C2"
$ws.Cells.Item(5, 2).Value  = "CL05A105KO5NNNC"
$ws.Cells.Item(5, 3).Value  = "C0402"
$ws.Cells.Item(5, 4).Value  = "18.544mm"
$ws.Cells.Item(5, 5).Value  = "-28.48mm"
$ws.Cells.Item(5, 6).Value  = "18.544mm"
$ws.Cells.Item(5, 7).Value  = "-28.48mm"
$ws.Cells.Item(5, 8).Value  = "19.09mm"
$ws.Cells.Item(5, 9).Value  = "-28.48mm"
$ws.Cells.Item(5, 10).Value = 2
$ws.Cells.Item(5, 11).Value = "T"
$ws.Cells.Item(5, 12).Value = 180
$ws.Cells.Item(5, 13).Value = "Yes"
$ws.Cells.Item(5, 14).Value = "1uF"

# --- Row 6: R1 resistor ---
$ws.Cells.Item(6, 1).Value  = "R1"
$ws.Cells.Item(6, 2).Value  = "RC0402FR-07100KL"
$ws.Cells.Item(6, 3).Value  = "R0402"
$ws.Cells.Item(6, 4).Value  = "25.4mm"
$ws.Cells.Item(6, 5).Value  = "-27.432mm"
$ws.Cells.Item(6, 6).Value  = "25.4mm"
$ws.Cells.Item(6, 7).Value  = "-27.432mm"
$ws.Cells.Item(6, 8).Value  = "25.4mm"
$ws.Cells.Item(6, 9).Value  = "-26.999mm"
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(6, 11).Value = "T"
$ws.Cells.Item(6, 12).Value = 270
$ws.Cells.Item(6, 13).Value = "Yes"
$ws.Cells.Item(6, 14).Value = "100kΩ"

# --- Row 7: U2 load switch ---
$ws.Cells.Item(7, 1).Value  = "U2"
$ws.Cells.Item(7, 2).Value  = "TPS22918DBVR"
$ws.Cells.Item(7, 3).Value  = "SOT-23-6_L2.9-W1.6-P0.95-LS2.8-BR"
$ws.Cells.Item(7, 4).Value  = "22.479mm"
$ws.Cells.Item(7, 5).Value  = "-29.464mm"
$ws.Cells.Item(7, 6).Value  = "22.479mm"
$ws.Cells.Item(7, 7).Value  = "-29.464mm"
$ws.Cells.Item(7, 8).Value  = "23.829mm"
$ws.Cells.Item(7, 9).Value  = "-30.414mm"
$ws.Cells.Item(7, 10).Value = 6
$ws.Cells.Item(7, 11).Value = "T"
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = "Yes"
$ws.Cells.Item(7, 14).Value = "TPS22918DBVR"
